# "excel data driven integration"
# Updates the two generated email addresses on the registerUsr sheet
# (the trailing "d" typo -> "y") and moves the sheet's active
# selection from G8 to N8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# registerUsr!G2 / G3 hold the test-data email addresses (still shown
# via the existing Email hyperlink style/format) - only the displayed
# text changes, the underlying mailto: hyperlinks are left as-is.
$ws.Range("G2").Value = "vijeysssdaay@gmail.com"
$ws.Range("G3").Value = "divyammmqqy@gmail.com"

# Move the sheet's selection from G8 to N8.
$ws.Range("N8").Select() | Out-Null
